$d = $word.ActiveDocument

# Insert the bulk of the new content (all new paragraphs) at the very start of the
# document body via a raw WordprocessingML package fragment.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:sz w:val="32"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="32"/>
    </w:rPr>
    <w:t>The (Failed) Creation of a Language Change</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="32"/>
    </w:rPr>
    <w:t xml:space="preserve"> Simulator</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t xml:space="preserve">The goal of this project is to create a software that simulates language change. Below is a list of steps I have to take towards completing the project. Of course, I was </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>wayyyy</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> too overly optimistic when I first envisioned the project…</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
    </w:rPr>
    <w:t>Describing the language</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>In order to simulate language change, we must first describe the current state of the language. Therefore, our software needs to be able to read in an</w:t>
  </w:r>
  <w:r>
    <w:t>d</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> st</w:t>
  </w:r>
  <w:r>
    <w:t>ore the grammar of the language. The user should be able to input descriptions of an existing language, or they could use the software to create their own!</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>1</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:tab/>
    <w:t>Phonemic inventory</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t xml:space="preserve">The most basic component of a (spoken) language is a phoneme – the smallest unit of sound. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Every language has a phonemic inventory, which is a collection of sounds that a speaker of that language perceives as distinctive. Each phoneme can be described using a feature matrix – a set of minimally necessary binary features that distinguishes it from all other phonemes in the language. </w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t xml:space="preserve">Thus, our software must be able to (1) store a phonemic inventory, (2) associate each phoneme to its set of features, and (3) allow the user to access phonemes via feature matrixes, and vice versa. As a bonus, we could also incorporate an algorithm that </w:t>
  </w:r>
  <w:r>
    <w:t>outputs a feature matrix given a list of feature values, and a system that allows the user to describe or create a phonetic-based writing system and associate each phoneme to an orthographical element.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>2</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:tab/>
    <w:t>Phonology</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:r>
    <w:t xml:space="preserve">Every language has phonological rules that govern how phonemes are realized in speech. </w:t>
  </w:r>
  <w:r>
    <w:t>These rules are encoded using the same feature matrixes described above.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>Thus, our software must be able to (1) allow the user to input and store phonological rules, and (2) correctly apply these rules to phonemic representations and yield the appropriate phonetic forms.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
    </w:rPr>
    <w:t>Time out…</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>So</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve"> by then I realized that even getting that far is a bit of a stretch… And I was right…</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
    </w:rPr>
    <w:t>An account of my unfortunate struggles</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>1</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:tab/>
    <w:t>Creating a user interface</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r = $d.Range(0, 0)
$r.InsertXML($xml)

# The original (only) paragraph -- which carries the _GoBack bookmark -- is now the
# last paragraph in the document. Add the new leading run to it in front of the
# bookmark, rather than introducing a brand new paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$insertPoint = $d.Range($lastRange.Start, $lastRange.Start)
$insertPoint.InsertBefore("When I first started the project, I decided that the first thing I needed to do was to create a user interface. ")

Write-Host "Paragraphs count:" $d.Paragraphs.Count
